$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Force text format on numeric-looking Price cells so exact formatting (trailing zeros, etc.) is preserved
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D51").NumberFormat = "@"

# Apply updated cell values
$ws.Range("D2").Value = '70.515.07'
$ws.Range("E2").Value = '  +4.38%  '
$ws.Range("D3").Value = '3.633.02'
$ws.Range("E3").Value = '  +3.96%  '
$ws.Range("E4").Value = '  +0.03%  '
$ws.Range("D5").Value = '593.64'
$ws.Range("E5").Value = '  +0.89%  '
$ws.Range("D6").Value = '195.24'
$ws.Range("E6").Value = '  +5.41%  '
$ws.Range("E7").Value = '  +1.78%  '
$ws.Range("D8").Value = '3.627.20'
$ws.Range("E8").Value = '  +3.85%  '
$ws.Range("E10").Value = '  +1.99%  '
$ws.Range("D11").Value = '0.671'
$ws.Range("E11").Value = '  +2.92%  '
$ws.Range("D12").Value = '58.66'
$ws.Range("E12").Value = '  +4.03%  '
$ws.Range("D13").Value = '0.0000292'
$ws.Range("E13").Value = '  +3.57%  '
$ws.Range("E14").Value = '  +4.93%  '
$ws.Range("D15").Value = '4.212.01'
$ws.Range("E15").Value = '  +4.24%  '
$ws.Range("D16").Value = '19.80'
$ws.Range("E16").Value = '  +5.29%  '
$ws.Range("D17").Value = '3.631.00'
$ws.Range("E17").Value = '  +4.01%  '
$ws.Range("D18").Value = '70.438.33'
$ws.Range("E18").Value = '  +4.34%  '
$ws.Range("D19").Value = '12.72'
$ws.Range("E19").Value = '  +4.17%  '
$ws.Range("E20").Value = '  +1.82%  '
$ws.Range("D21").Value = '1.07'
$ws.Range("E21").Value = '  +4.31%  '
$ws.Range("D22").Value = '489.37'
$ws.Range("E22").Value = '  -0.67%  '
$ws.Range("D23").Value = '19.40'
$ws.Range("E23").Value = '  +16.20%  '
$ws.Range("D24").Value = '5.39'
$ws.Range("E24").Value = '  -1.10%  '
$ws.Range("E25").Value = '  -0.40%  '
$ws.Range("D26").Value = '91.27'
$ws.Range("E26").Value = '  +1.00%  '
$ws.Range("D27").Value = '3.16'
$ws.Range("E27").Value = '  +6.07%  '
$ws.Range("D28").Value = '11.48'
$ws.Range("E28").Value = '  +4.22%  '
$ws.Range("E29").Value = '  +5.50%  '
$ws.Range("D30").Value = '33.04'
$ws.Range("E30").Value = '  +4.28%  '
$ws.Range("D31").Value = '7.91'
$ws.Range("E31").Value = '  +10.26%  '
$ws.Range("D32").Value = '628.47'
$ws.Range("E32").Value = '  +5.43%  '
$ws.Range("D33").Value = '0.121'
$ws.Range("E33").Value = '  +7.87%  '
$ws.Range("D34").Value = '12.31'
$ws.Range("E34").Value = '  +4.47%  '
$ws.Range("D35").Value = '66.19'
$ws.Range("E35").Value = '  +2.56%  '
$ws.Range("D36").Value = '40.57'
$ws.Range("E36").Value = '  +10.54%  '
$ws.Range("D37").Value = '0.415'
$ws.Range("E37").Value = '  +6.72%  '
$ws.Range("D38").Value = '0.0₃0824'
$ws.Range("E38").Value = '  +6.78%  '
$ws.Range("E39").Value = '  -0.02%  '
$ws.Range("D40").Value = '0.147'
$ws.Range("E40").Value = '  -2.26%  '
$ws.Range("D41").Value = '3.61'
$ws.Range("E41").Value = '  +1.19%  '
$ws.Range("D42").Value = '3.294.14'
$ws.Range("E42").Value = '  +1.01%  '
$ws.Range("D43").Value = '3.15'
$ws.Range("E43").Value = '  +7.68%  '
$ws.Range("D44").Value = '2.83'
$ws.Range("E44").Value = '  +11.38%  '
$ws.Range("E45").Value = '  +5.23%  '
$ws.Range("B46").Value = 'dogwifhat'
$ws.Range("C46").Value = 'https://coinranking.com/coin/sZUrmToWF+dogwifhat-wif'
$ws.Range("D46").Value = '2.83'
$ws.Range("E46").Value = '  +2.92%  '
$ws.Range("B47").Value = 'ApeXProtocol'
$ws.Range("C47").Value = 'https://coinranking.com/coin/ze0N2Rcyu+apexprotocol-apex'
$ws.Range("D47").Value = '3.31'
$ws.Range("E47").Value = '  +0.30%  '
$ws.Range("B48").Value = 'Stellar'
$ws.Range("C48").Value = 'https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm'
$ws.Range("D48").Value = '0.139'
$ws.Range("E48").Value = '  +2.52%  '
$ws.Range("D49").Value = '9.21'
$ws.Range("E49").Value = '  +5.11%  '
$ws.Range("E50").Value = '  +3.42%  '
$ws.Range("B51").Value = 'Monero'
$ws.Range("C51").Value = 'https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr'
$ws.Range("D51").Value = '143.31'
$ws.Range("E51").Value = '  +1.68%  '
